$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "Andrea"
$ws.Range("G6").Value = 28
$ws.Range("F7").Value = "Samuel"
$ws.Range("G7").Value = 40

$ws.Range("F8").Select()
